$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two email/contact cells in column F with new placeholder text
$ws.Range("F1").Value = "Seu Madruga: email_Generico"
$ws.Range("F2").Value = "Chavez:  email_Generico"

# Move the active selection from A2 to F4
$ws.Range("F4").Select()
